$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.75
$ws.Range("I2").Value = 2.9
$ws.Range("T2").Value = 6
$ws.Range("AE2").Value = 12

# Row 8
$ws.Range("N8").Value = 1.95
$ws.Range("O8").Value = 1.85

# Row 14
$ws.Range("G14").Value = 1.78
$ws.Range("I14").Value = 4.5
$ws.Range("T14").Value = 4.9
$ws.Range("U14").Value = 6.8
$ws.Range("W14").Value = 13.5
$ws.Range("Z14").Value = 6.5
$ws.Range("AA14").Value = 6.7
$ws.Range("AD14").Value = 9
$ws.Range("AE14").Value = 23
$ws.Range("AF14").Value = 17
$ws.Range("AI14").Value = 90

# Row 15
$ws.Range("G15").Value = 2.52
$ws.Range("I15").Value = 2.77
$ws.Range("T15").Value = 6
$ws.Range("U15").Value = 10.75
$ws.Range("W15").Value = 27
$ws.Range("X15").Value = 28
$ws.Range("Y15").Value = 55
$ws.Range("AA15").Value = 6
$ws.Range("AB15").Value = 20
$ws.Range("AE15").Value = 12
$ws.Range("AF15").Value = 11.5
$ws.Range("AG15").Value = 32

# Row 17
$ws.Range("G17").Value = 2.63
$ws.Range("I17").Value = 2.63
$ws.Range("K17").Value = 8
$ws.Range("T17").Value = 7
$ws.Range("Y17").Value = 41
$ws.Range("AC17").Value = 67
$ws.Range("AD17").Value = 7
$ws.Range("AG17").Value = 26
